# Update averaged-questions-arm-study-data sheet with corrected values.
# The "_Competency" and "_Discomfort" rows for each study id had their
# labels/values swapped and the underlying averages recomputed; this
# script writes the corrected label + B:E averages for rows 2-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "1_Discomfort", 1.666666666666667, 4.166666666666667, 4.333333333333333, 5.5),
    @(3, "1_Competency", 6.833333333333333, 6.833333333333333, 1.666666666666667, 4.5),
    @(4, "1_Safety", 3, 3.666666666666667, 2.666666666666667, 2.333333333333333),
    @(5, "2_Discomfort", 5.166666666666667, 6.833333333333333, 4.333333333333333, 5.833333333333333),
    @(6, "2_Competency", 6, 3.833333333333333, 5.5, 5.833333333333333),
    @(7, "2_Safety", 3, 3.666666666666667, 2.333333333333333, 3.333333333333333),
    @(8, "3_Discomfort", 1.5, 3.666666666666667, 5.833333333333333, 6.666666666666667),
    @(9, "3_Competency", 7.166666666666667, 6.166666666666667, 1.833333333333333, 2.666666666666667),
    @(10, "3_Safety", 3.666666666666667, 2.333333333333333, 3.333333333333333, 2.333333333333333),
    @(11, "4_Discomfort", 1.5, 6.666666666666667, 4.666666666666667, 3.5),
    @(12, "4_Competency", 7, 2.666666666666667, 1.666666666666667, 4.5),
    @(13, "4_Safety", 3.333333333333333, 3.666666666666667, 2.666666666666667, 2.333333333333333),
    @(14, "5_Discomfort", 3, 7.166666666666667, 3.5, 2.5),
    @(15, "5_Competency", 6.833333333333333, 4.166666666666667, 2.166666666666667, 5.166666666666667),
    @(16, "5_Safety", 2, 3.666666666666667, 2.666666666666667, 2.666666666666667),
    @(17, "6_Discomfort", 3.833333333333333, 6.166666666666667, 3.5, 3.333333333333333),
    @(18, "6_Competency", 6.666666666666667, 2.166666666666667, 5.5, 4.333333333333333),
    @(19, "6_Safety", 3.666666666666667, 3.666666666666667, 2.666666666666667, 2.666666666666667),
    @(20, "7_Discomfort", 3, 2.333333333333333, 2.333333333333333, 1.333333333333333),
    @(21, "7_Competency", 7, 6, 6.666666666666667, 6.5),
    @(22, "7_Safety", 3, 2.333333333333333, 2.333333333333333, 2.333333333333333),
    @(23, "8_Discomfort", 2.333333333333333, 4, 3.5, 4.666666666666667),
    @(24, "8_Competency", 7.5, 7, 4.833333333333333, 7.333333333333333),
    @(25, "8_Safety", 3, 2.333333333333333, 2, 3),
    @(26, "9_Discomfort", 2.833333333333333, 2.833333333333333, 3.333333333333333, 2),
    @(27, "9_Competency", 7, 6.666666666666667, 5.666666666666667, 7.333333333333333),
    @(28, "9_Safety", 3, 2.333333333333333, 2.333333333333333, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
